# Regenerate save_data to use K instead of Strike# (calc and write s_vals)
# Applies updated strikeout (K) counts in column G of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 2
    9  = 2
    10 = 0
    11 = 1
    12 = 1
    13 = 0
    14 = 2
    15 = 1
    16 = 2
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 2
    25 = 0
    26 = 1
    27 = 1
    28 = 2
    29 = 3
    30 = 2
    31 = 1
    32 = 2
    34 = 3
    35 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
